# Fruta / hortaliza, semanal
# Insert 2 new weekly price rows for Ciruela "Angeleno" (Primera/Segunda) dated 2023-03-03
# (serial 44988) at the top of the existing block (row 278), pushing the rest of the
# table down by two rows (old row 355 -> new row 357).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 278; this shifts rows 278:355 -> 280:357
# and copies formatting (incl. the date-style on column D) from the row above, matching
# the original workbook's row layout.
$ws.Rows("278:279").Insert()

# New row 278: Ciruela - Angeleno - Primera
$ws.Range("A278").Value = 8
$ws.Range("B278").Value = "Terminal La Palmera de La Serena"
$ws.Range("C278").Value = "Coquimbo"
$ws.Range("D278").Value = 44988
$ws.Range("E278").Value = 4
$ws.Range("F278").Value = "Fruta"
$ws.Range("G278").Value = 100103
$ws.Range("H278").Value = "Frutos de hueso (carozo)"
$ws.Range("I278").Value = 100103002
$ws.Range("J278").Value = "Ciruela"
$ws.Range("K278").Value = "Angeleno"
$ws.Range("L278").Value = "Primera"
$ws.Range("M278").Value = 16
$ws.Range("N278").Value = 170000
$ws.Range("O278").Value = 180000
$ws.Range("P278").Value = 175000
$ws.Range("Q278").Value = "$/bins (450 kilos)"
$ws.Range("R278").Value = "Región de O'Higgins"
$ws.Range("S278").Value = 389
$ws.Range("T278").Value = 450

# New row 279: Ciruela - Angeleno - Segunda
$ws.Range("A279").Value = 8
$ws.Range("B279").Value = "Terminal La Palmera de La Serena"
$ws.Range("C279").Value = "Coquimbo"
$ws.Range("D279").Value = 44988
$ws.Range("E279").Value = 4
$ws.Range("F279").Value = "Fruta"
$ws.Range("G279").Value = 100103
$ws.Range("H279").Value = "Frutos de hueso (carozo)"
$ws.Range("I279").Value = 100103002
$ws.Range("J279").Value = "Ciruela"
$ws.Range("K279").Value = "Angeleno"
$ws.Range("L279").Value = "Segunda"
$ws.Range("M279").Value = 20
$ws.Range("N279").Value = 140000
$ws.Range("O279").Value = 150000
$ws.Range("P279").Value = 145000
$ws.Range("Q279").Value = "$/bins (450 kilos)"
$ws.Range("R279").Value = "Región de O'Higgins"
$ws.Range("S279").Value = 322
$ws.Range("T279").Value = 450
